# Auto-generated Excel COM-interop script
# Applies per-cell updates to ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets
# matching the scheduled-runner profit recompute diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1823.5333
$ws.Range("I4").Value = 716.7778
$ws.Range("J4").Value = 3483.6667
$ws.Range("K4").Value = 716.7778
$ws.Range("L4").Value = 3483.6667
$ws.Range("M4").Value = -602.7778
$ws.Range("N4").Value = -3711.6667

$ws.Range("H39").Value = 499.5
$ws.Range("I39").Value = 398
$ws.Range("J39").Value = 533.3333
$ws.Range("K39").Value = 1194
$ws.Range("L39").Value = 1599.9999
$ws.Range("M39").Value = -898
$ws.Range("N39").Value = -2191.9999

$ws.Range("H62").Value = 6644.4375
$ws.Range("I62").Value = 3067.0833
$ws.Range("J62").Value = 17376.5
$ws.Range("K62").Value = 3067.0833
$ws.Range("L62").Value = 17376.5
$ws.Range("M62").Value = -2443.0833
$ws.Range("N62").Value = -18624.5

$ws.Range("H65").Value = 6644.4375
$ws.Range("I65").Value = 3067.0833
$ws.Range("J65").Value = 17376.5
$ws.Range("K65").Value = 15335.4165
$ws.Range("L65").Value = 86882.5
$ws.Range("M65").Value = -12215.4165
$ws.Range("N65").Value = -93122.5

$ws.Range("H141").Value = 3027.1724
$ws.Range("I141").Value = 2100.5264
$ws.Range("J141").Value = 4787.8
$ws.Range("K141").Value = 6301.5792
$ws.Range("L141").Value = 14363.4
$ws.Range("M141").Value = -1121.5792
$ws.Range("N141").Value = -24723.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 31933.334
$ws.Range("J24").Value = 31933.334
$ws.Range("L24").Value = 31933.334
$ws.Range("N24").Value = -32681.334

$ws.Range("H100").Value = 31933.334
$ws.Range("J100").Value = 31933.334
$ws.Range("L100").Value = 31933.334
$ws.Range("N100").Value = -34097.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4474.109
$ws.Range("I31").Value = 3985.4285
$ws.Range("J31").Value = 6029
$ws.Range("K31").Value = 3985.4285
$ws.Range("L31").Value = 6029
$ws.Range("M31").Value = -3690.4285
$ws.Range("N31").Value = -6619

$ws.Range("H34").Value = 4474.109
$ws.Range("I34").Value = 3985.4285
$ws.Range("J34").Value = 6029
$ws.Range("K34").Value = 3985.4285
$ws.Range("L34").Value = 6029
$ws.Range("M34").Value = -3783.4285
$ws.Range("N34").Value = -6433

$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793

$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373

$ws.Range("H74").Value = 36656.75
$ws.Range("J74").Value = 36656.75
$ws.Range("L74").Value = 36656.75
$ws.Range("N74").Value = -38404.75

$ws.Range("H77").Value = 36656.75
$ws.Range("J77").Value = 36656.75
$ws.Range("L77").Value = 109970.25
$ws.Range("N77").Value = -118706.25

$ws.Range("H99").Value = 3585.5
$ws.Range("I99").Value = 3399.5
$ws.Range("J99").Value = 3678.5
$ws.Range("K99").Value = 3399.5
$ws.Range("L99").Value = 3678.5
$ws.Range("M99").Value = -1901.5
$ws.Range("N99").Value = -6674.5

$ws.Range("H106").Value = 38000
$ws.Range("J106").Value = 38000
$ws.Range("L106").Value = 38000
$ws.Range("N106").Value = -40524

$ws.Range("H122").Value = 9972.157999999999
$ws.Range("I122").Value = 4657.1333
$ws.Range("J122").Value = 29903.5
$ws.Range("K122").Value = 13971.3999
$ws.Range("L122").Value = 89710.5
$ws.Range("M122").Value = -11521.3999
$ws.Range("N122").Value = -94610.5

$ws.Range("H126").Value = 3585.5
$ws.Range("I126").Value = 3399.5
$ws.Range("J126").Value = 3678.5
$ws.Range("K126").Value = 10198.5
$ws.Range("L126").Value = 11035.5
$ws.Range("M126").Value = -7728.5
$ws.Range("N126").Value = -15975.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6628
$ws.Range("J5").Value = 25356.25
$ws.Range("L5").Value = 76068.75
$ws.Range("N5").Value = -76292.75

$ws.Range("H135").Value = 6628
$ws.Range("J135").Value = 25356.25
$ws.Range("L135").Value = 228206.25
$ws.Range("N135").Value = -233276.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 28900
$ws.Range("J46").Value = 28900
$ws.Range("L46").Value = 28900
$ws.Range("N46").Value = -29212

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6938.125
$ws.Range("I7").Value = 6500
$ws.Range("J7").Value = 7668.3335
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7668.3335
$ws.Range("M7").Value = -6388
$ws.Range("N7").Value = -7892.3335

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H35").Value = 20891.584
$ws.Range("I35").Value = 9979.799999999999
$ws.Range("J35").Value = 28685.715
$ws.Range("K35").Value = 9979.799999999999
$ws.Range("L35").Value = 28685.715
$ws.Range("M35").Value = -9643.799999999999
$ws.Range("N35").Value = -29357.715

$ws.Range("H43").Value = 19980
$ws.Range("J43").Value = 19980
$ws.Range("L43").Value = 19980
$ws.Range("N43").Value = -20366

$ws.Range("H104").Value = 17666.666
$ws.Range("J104").Value = 17666.666
$ws.Range("L104").Value = 17666.666
$ws.Range("N104").Value = -24654.666

$ws.Range("H126").Value = 6938.125
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 7668.3335
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 23005.0005
$ws.Range("M126").Value = -17030
$ws.Range("N126").Value = -27945.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 6200
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 6600
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 6600
$ws.Range("M29").Value = -4710
$ws.Range("N29").Value = -7180

$ws.Range("H101").Value = 25777.6
$ws.Range("J101").Value = 25777.6
$ws.Range("L101").Value = 25777.6
$ws.Range("N101").Value = -32267.6

$ws.Range("H104").Value = 20790
$ws.Range("J104").Value = 20790
$ws.Range("L104").Value = 20790
$ws.Range("N104").Value = -27778

$ws.Range("H122").Value = 4365.303
$ws.Range("I122").Value = 1266.32
$ws.Range("J122").Value = 14049.625
$ws.Range("K122").Value = 3798.96
$ws.Range("L122").Value = 42148.875
$ws.Range("M122").Value = -1348.96
$ws.Range("N122").Value = -47048.875

$ws.Range("H126").Value = 1217.1923
$ws.Range("I126").Value = 1170.8948
$ws.Range("J126").Value = 1342.8572
$ws.Range("K126").Value = 3512.6844
$ws.Range("L126").Value = 4028.5716
$ws.Range("M126").Value = -1042.6844
$ws.Range("N126").Value = -8968.571599999999
